# Display - Raw.xlsx : refreshed raw measurement data (mean/median increase)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw values for A2:A31 (column A holds the "Standard" display measurements)
$values = @(
    168.11712,
    185.24016,
    182.12688,
    179.0136,
    183.68351999999999,
    171.2304,
    180.57023999999899,
    180.57023999999899,
    180.57023999999899,
    179.0136,
    180.57023999999899,
    179.0136,
    177.45696000000001,
    177.45696000000001,
    179.0136,
    177.45696000000001,
    179.0136,
    179.0136,
    182.12688,
    177.45696000000001,
    180.57023999999899,
    171.2304,
    179.0136,
    180.57023999999899,
    180.57023999999899,
    180.57023999999899,
    180.57023999999899,
    182.12688,
    179.0136,
    179.0136
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Match the saved cursor/selection position recorded after the refresh
$ws.Range("D27").Select()

$wb.Save()
